$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.954.55"
$ws.Range("E2").Value = "  +2.01%  "

$ws.Range("D3").Value = "1.814.18"
$ws.Range("E3").Value = "  +2.50%  "

$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").Value = "'312.18"
$ws.Range("E5").Value = "  +1.88%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").Value = "'0.4296"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "'0.3672"
$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("D9").Value = "'0.07231"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").Value = "2.232.04"
$ws.Range("E10").Value = "  +24.65%  "

$ws.Range("D11").Value = "'0.8628"
$ws.Range("E11").Value = "  +1.52%  "

$ws.Range("D12").Value = "'21.21"
$ws.Range("E12").Value = "  +4.13%  "

$ws.Range("D13").Value = "'5.398"
$ws.Range("E13").Value = "  +3.12%  "

$ws.Range("D14").Value = "'6.597"
$ws.Range("E14").Value = "  +2.59%  "

$ws.Range("D15").Value = "'0.06946"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").Value = "'81.14"
$ws.Range("E16").Value = "  +2.30%  "

$ws.Range("D17").Value = "'1.011"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").Value = "'0.000008881"
$ws.Range("E18").Value = "  +2.33%  "

$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("D20").Value = "'15.18"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").Value = "27.000.67"
$ws.Range("E21").Value = "  +2.13%  "

$ws.Range("E22").Value = "  +1.30%  "

$ws.Range("D23").Value = "2.450.67"
$ws.Range("E23").Value = "  +22.12%  "

$ws.Range("D24").Value = "'11.00"
$ws.Range("E24").Value = "  -1.98%  "

$ws.Range("D25").Value = "'153.87"
$ws.Range("E25").Value = "  +1.24%  "

$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("E27").Value = "  +1.15%  "

$ws.Range("D28").Value = "'5.222"
$ws.Range("E28").Value = "  +2.71%  "

$ws.Range("D29").Value = "'1.902"
$ws.Range("E29").Value = "  +8.64%  "

$ws.Range("D30").Value = "'114.55"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("D31").Value = "'0.08945"
$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("E32").Value = "  +7.07%  "

$ws.Range("D33").Value = "'0.7485"
$ws.Range("E33").Value = "  +3.25%  "

$ws.Range("D34").Value = "'4.419"
$ws.Range("E34").Value = "  +1.98%  "

$ws.Range("D35").Value = "'2.805"
$ws.Range("E35").Value = "  +2.17%  "

$ws.Range("E36").Value = "  +0.26%  "

$ws.Range("E37").Value = "  +4.03%  "

$ws.Range("D38").Value = "'0.05208"

$ws.Range("D39").Value = "'0.01920"
$ws.Range("E39").Value = "  +1.66%  "

$ws.Range("D40").Value = "'0.5103"
$ws.Range("E40").Value = "  +3.67%  "

$ws.Range("D41").Value = "'2.749"
$ws.Range("E41").Value = "  +6.89%  "

$ws.Range("D42").Value = "'0.1652"
$ws.Range("E42").Value = "  +2.96%  "

$ws.Range("D43").Value = "'6.473"
$ws.Range("E43").Value = "  +3.45%  "

$ws.Range("D44").Value = "'8.328"
$ws.Range("E44").Value = "  +4.05%  "

$ws.Range("D45").Value = "'106.77"
$ws.Range("E45").Value = "  +1.81%  "

$ws.Range("E46").Value = "  +1.94%  "

$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("D48").Value = "'0.4581"
$ws.Range("E48").Value = "  +1.95%  "

$ws.Range("E49").Value = "  +3.31%  "

$ws.Range("D50").Value = "'0.06213"
$ws.Range("E50").Value = "  +0.32%  "

$ws.Range("D51").Value = "'1.850"
$ws.Range("E51").Value = "  +6.47%  "
